# Update market-price / profit figures (columns H:N) on several Leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, refreshed by the
# scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")

# Row 43 - Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 80 - Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 3569.8823
$ws.Range("I80").Value = 2697.3333
$ws.Range("J80").Value = 4045.818
$ws.Range("K80").Value = 8091.999899999999
$ws.Range("L80").Value = 12137.454
$ws.Range("M80").Value = -7093.999899999999
$ws.Range("N80").Value = -14133.454

# Row 83 - Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 3569.8823
$ws.Range("I83").Value = 2697.3333
$ws.Range("J83").Value = 4045.818
$ws.Range("K83").Value = 24275.9997
$ws.Range("L83").Value = 36412.362
$ws.Range("M83").Value = -19283.9997
$ws.Range("N83").Value = -46396.362

# Row 125 - Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 285730720
$ws.Range("I125").Value = 375014500
$ws.Range("J125").Value = 166685660
$ws.Range("K125").Value = 3375130500
$ws.Range("L125").Value = 1500170940
$ws.Range("M125").Value = -3375128040
$ws.Range("N125").Value = -1500175860

# Row 132 - Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 492.0625
$ws.Range("I132").Value = 498.2
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 1494.6
$ws.Range("L132").Value = 1200
$ws.Range("M132").Value = 1035.4
$ws.Range("N132").Value = -6260

# Row 135 - For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 25
$ws.Range("I135").Value = 25
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 225
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = 2310

# Row 137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1590
$ws.Range("I137").Value = 1070
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 3210
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -660
$ws.Range("N137").Value = -12600

# Row 138 - All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2199.2273
$ws.Range("I138").Value = 2093
$ws.Range("J138").Value = 2352.6667
$ws.Range("K138").Value = 6279
$ws.Range("L138").Value = 7058.000100000001
$ws.Range("M138").Value = -1139
$ws.Range("N138").Value = -17338.0001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")

# Row 5 - The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 2552
$ws.Range("I5").Value = 2552
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2552
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2440
$ws.Range("N5").ClearContents()

# Row 61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 5240.8184
$ws.Range("I61").Value = 3562
$ws.Range("J61").Value = 7665.778
$ws.Range("K61").Value = 3562
$ws.Range("L61").Value = 7665.778
$ws.Range("M61").Value = -3350
$ws.Range("N61").Value = -8089.778

# Row 74 - As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2994.75
$ws.Range("I74").Value = 2993
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2993
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2119
$ws.Range("N74").Value = -4748

# Row 77 - Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2994.75
$ws.Range("I77").Value = 2993
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 14965
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -10597
$ws.Range("N77").Value = -23736

# Row 136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 5240.8184
$ws.Range("I136").Value = 3562
$ws.Range("J136").Value = 7665.778
$ws.Range("K136").Value = 10686
$ws.Range("L136").Value = 22997.334
$ws.Range("M136").Value = -8136
$ws.Range("N136").Value = -28097.334

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")

# Row 4 - Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 2552
$ws.Range("I4").Value = 2552
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2552
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2437
$ws.Range("N4").ClearContents()

# Row 134 - Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1651.6
$ws.Range("I134").Value = 1651.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4954.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2419.799999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1537.4445
$ws.Range("I31").Value = 1504.625
$ws.Range("J31").Value = 1800
$ws.Range("K31").Value = 1504.625
$ws.Range("L31").Value = 1800
$ws.Range("M31").Value = -1209.625
$ws.Range("N31").Value = -2390

# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1537.4445
$ws.Range("I34").Value = 1504.625
$ws.Range("J34").Value = 1800
$ws.Range("K34").Value = 1504.625
$ws.Range("L34").Value = 1800
$ws.Range("M34").Value = -1302.625
$ws.Range("N34").Value = -2204

# Row 132 - Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2938.4
$ws.Range("I132").Value = 3198
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 9594
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -7064
$ws.Range("N132").Value = -10760

# Row 134 - Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 634.7059
$ws.Range("I134").Value = 634.7059
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1904.1177
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 630.8822999999998
$ws.Range("N134").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")

# Row 68 - Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 12112
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 12112
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 36336
$ws.Range("N68").Value = -37958

# Row 71 - No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 12112
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 12112
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 109008
$ws.Range("N71").Value = -117120

# Row 107 - Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 485.66666
$ws.Range("I107").Value = 250
$ws.Range("J107").Value = 553
$ws.Range("K107").Value = 750
$ws.Range("L107").Value = 1659
$ws.Range("M107").Value = 1170
$ws.Range("N107").Value = -5499

# Row 116 - On a Full Stomach / Sausage Links
$ws.Range("H116").Value = 2964.5
$ws.Range("I116").Value = 2964.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 8893.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -5451.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")

# Row 80 - Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 1488.5555
$ws.Range("I80").Value = 1174.625
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 1174.625
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -176.625
$ws.Range("N80").Value = -5996

# Row 83 - With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 1488.5555
$ws.Range("I83").Value = 1174.625
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 5873.125
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -881.125
$ws.Range("N83").Value = -29984

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")

# Row 16 - Saddle Sore / Hard Leather
$ws.Range("H16").Value = 1786
$ws.Range("I16").Value = 1826.8572
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1826.8572
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1656.8572
$ws.Range("N16").Value = -1840

# Row 46 - Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 33992.312
$ws.Range("I46").Value = 64434.5
$ws.Range("J46").Value = 3550.125
$ws.Range("K46").Value = 64434.5
$ws.Range("L46").Value = 3550.125
$ws.Range("M46").Value = -64246.5
$ws.Range("N46").Value = -3926.125

# Row 132 - Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 1527.9333
$ws.Range("I132").Value = 1301.28
$ws.Range("J132").Value = 2661.2
$ws.Range("K132").Value = 3903.84
$ws.Range("L132").Value = 7983.599999999999
$ws.Range("M132").Value = -1373.84
$ws.Range("N132").Value = -13043.6

# Row 136 - Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4620.4707
$ws.Range("I136").Value = 4349.615
$ws.Range("J136").Value = 5500.75
$ws.Range("K136").Value = 13048.845
$ws.Range("L136").Value = 16502.25
$ws.Range("M136").Value = -10498.845
$ws.Range("N136").Value = -21602.25

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")

# Row 62 - Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 12605.214
$ws.Range("I62").Value = 12111.333
$ws.Range("J62").Value = 13494.2
$ws.Range("K62").Value = 12111.333
$ws.Range("L62").Value = 13494.2
$ws.Range("M62").Value = -11487.333
$ws.Range("N62").Value = -14742.2

# Row 65 - Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 12605.214
$ws.Range("I65").Value = 12111.333
$ws.Range("J65").Value = 13494.2
$ws.Range("K65").Value = 60556.665
$ws.Range("L65").Value = 67471
$ws.Range("M65").Value = -57436.665
$ws.Range("N65").Value = -73711

# Row 132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2367.3547
$ws.Range("I132").Value = 2521.8076
$ws.Range("J132").Value = 1564.2
$ws.Range("K132").Value = 7565.4228
$ws.Range("L132").Value = 4692.6
$ws.Range("M132").Value = -5035.4228
$ws.Range("N132").Value = -9752.6

# Row 136 - Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 3664.1177
$ws.Range("I136").Value = 3664.1177
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8442.3531
$ws.Range("N136").ClearContents()
